$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H38").Value = 2150849.5
$ws.Range("I38").Value = 2481474
$ws.Range("J38").Value = 1790
$ws.Range("K38").Value = 7444422
$ws.Range("L38").Value = 5370
$ws.Range("M38").Value = -7444050
$ws.Range("N38").Value = -6114
$ws.Range("H58").Value = 6646250
$ws.Range("I58").Value = 3921938.5
$ws.Range("J58").Value = 7694062
$ws.Range("K58").Value = 11765815.5
$ws.Range("L58").Value = 23082186
$ws.Range("M58").Value = -11765665.5
$ws.Range("N58").Value = -23082486
$ws.Range("H62").Value = 2424.158
$ws.Range("I62").Value = 2316.1875
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2316.1875
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1692.1875
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2424.158
$ws.Range("I65").Value = 2316.1875
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 11580.9375
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -8460.9375
$ws.Range("N65").Value = -21240
$ws.Range("H88").Value = 2114.4092
$ws.Range("I88").Value = 685.375
$ws.Range("J88").Value = 2931
$ws.Range("K88").Value = 685.375
$ws.Range("L88").Value = 2931
$ws.Range("M88").Value = -279.375
$ws.Range("N88").Value = -3743
$ws.Range("H91").Value = 2114.4092
$ws.Range("I91").Value = 685.375
$ws.Range("J91").Value = 2931
$ws.Range("K91").Value = 685.375
$ws.Range("L91").Value = 2931
$ws.Range("M91").Value = 718.625
$ws.Range("N91").Value = -5739
$ws.Range("H123").Value = 45417
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 45417
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 45417
$ws.Range("N123").Value = -55217
$ws.Range("H135").Value = 3135.0244
$ws.Range("I135").Value = 727.3913
$ws.Range("J135").Value = 6211.4443
$ws.Range("K135").Value = 6546.5217
$ws.Range("L135").Value = 55902.9987
$ws.Range("M135").Value = -4011.5217
$ws.Range("N135").Value = -60972.9987
$ws.Range("H138").Value = 2886.8877
$ws.Range("I138").Value = 1498.2122
$ws.Range("J138").Value = 3591.9077
$ws.Range("K138").Value = 4494.6366
$ws.Range("L138").Value = 10775.7231
$ws.Range("M138").Value = 645.3634000000002
$ws.Range("N138").Value = -21055.7231

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 25538.357
$ws.Range("I32").Value = 9451.638
$ws.Range("J32").Value = 111334.2
$ws.Range("K32").Value = 9451.638
$ws.Range("L32").Value = 111334.2
$ws.Range("M32").Value = -9164.638
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H44").Value = 5700
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 5700
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 5700
$ws.Range("N44").Value = -6676
$ws.Range("H88").Value = 1750
$ws.Range("I88").Value = 1750
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1750
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1344
$ws.Range("H91").Value = 1750
$ws.Range("I91").Value = 1750
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1750
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -346
$ws.Range("H95").Value = 23984.666
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 23984.666
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 23984.666
$ws.Range("N95").Value = -29476.666
$ws.Range("H135").Value = 45074.145
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 45074.145
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 45074.145
$ws.Range("N135").Value = -55214.145
$ws.Range("H140").Value = 60825.445
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 60825.445
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 60825.445
$ws.Range("N140").Value = -71185.445

$ws = $wb.Worksheets("BSM")
$ws.Range("H35").Value = 19046.8
$ws.Range("I35").Value = 17000
$ws.Range("J35").Value = 19558.5
$ws.Range("K35").Value = 17000
$ws.Range("L35").Value = 19558.5
$ws.Range("M35").Value = -16690
$ws.Range("N35").Value = -20178.5
$ws.Range("H80").Value = 1089.2122
$ws.Range("I80").Value = 679.9231
$ws.Range("J80").Value = 1355.25
$ws.Range("K80").Value = 679.9231
$ws.Range("L80").Value = 1355.25
$ws.Range("M80").Value = 318.0769
$ws.Range("N80").Value = -3351.25
$ws.Range("H83").Value = 1089.2122
$ws.Range("I83").Value = 679.9231
$ws.Range("J83").Value = 1355.25
$ws.Range("K83").Value = 3399.6155
$ws.Range("L83").Value = 6776.25
$ws.Range("M83").Value = 1592.3845
$ws.Range("N83").Value = -16760.25
$ws.Range("H94").Value = 850
$ws.Range("I94").Value = 500
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 500
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -49
$ws.Range("N94").Value = -2102
$ws.Range("H105").Value = 401400.2
$ws.Range("I105").Value = 334996.66
$ws.Range("J105").Value = 501005.5
$ws.Range("K105").Value = 334996.66
$ws.Range("L105").Value = 501005.5
$ws.Range("M105").Value = -333249.66
$ws.Range("N105").Value = -504499.5
$ws.Range("H107").Value = 83334700
$ws.Range("I107").Value = 125001310
$ws.Range("J107").Value = 1485
$ws.Range("K107").Value = 125001310
$ws.Range("L107").Value = 1485
$ws.Range("M107").Value = -124999390
$ws.Range("N107").Value = -5325
$ws.Range("H137").Value = 37109.668
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 37109.668
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 37109.668
$ws.Range("N137").Value = -47309.668

$ws = $wb.Worksheets("CRP")
$ws.Range("H58").Value = 2121.6365
$ws.Range("I58").Value = 1980.125
$ws.Range("J58").Value = 2499
$ws.Range("K58").Value = 1980.125
$ws.Range("L58").Value = 2499
$ws.Range("M58").Value = -1777.125
$ws.Range("H134").Value = 1845.6364
$ws.Range("I134").Value = 979.8
$ws.Range("J134").Value = 2567.1667
$ws.Range("K134").Value = 2939.4
$ws.Range("L134").Value = 7701.500100000001
$ws.Range("M134").Value = -404.3999999999996
$ws.Range("N134").Value = -12771.5001
$ws.Range("H136").Value = 2121.6365
$ws.Range("I136").Value = 1980.125
$ws.Range("J136").Value = 2499
$ws.Range("K136").Value = 5940.375
$ws.Range("L136").Value = 7497
$ws.Range("M136").Value = -3390.375

$ws = $wb.Worksheets("CUL")
$ws.Range("H98").Value = 92006.18
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 92006.18
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 276018.54
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -279014.54
$ws.Range("H113").Value = 879.6129
$ws.Range("I113").Value = 1096.9412
$ws.Range("J113").Value = 615.7143
$ws.Range("K113").Value = 3290.8236
$ws.Range("L113").Value = 1847.1429
$ws.Range("M113").Value = -1120.8236
$ws.Range("N113").Value = -6187.1429
$ws.Range("H131").Value = 835.32
$ws.Range("I131").Value = 482.83334
$ws.Range("J131").Value = 912.6951
$ws.Range("K131").Value = 1448.50002
$ws.Range("L131").Value = 2738.0853
$ws.Range("M131").Value = 3591.49998
$ws.Range("N131").Value = -12818.0853

$ws = $wb.Worksheets("GSM")
$ws.Range("H70").Value = 161185.16
$ws.Range("I70").Value = 256299.75
$ws.Range("J70").Value = 9001.8
$ws.Range("K70").Value = 256299.75
$ws.Range("L70").Value = 9001.8
$ws.Range("M70").Value = -256029.75
$ws.Range("N70").Value = -9541.8
$ws.Range("H73").Value = 161185.16
$ws.Range("I73").Value = 256299.75
$ws.Range("J73").Value = 9001.8
$ws.Range("K73").Value = 256299.75
$ws.Range("L73").Value = 9001.8
$ws.Range("M73").Value = -255363.75
$ws.Range("N73").Value = -10873.8
$ws.Range("H107").Value = 673830.94
$ws.Range("I107").Value = 357.45456
$ws.Range("J107").Value = 2525883
$ws.Range("K107").Value = 357.45456
$ws.Range("L107").Value = 2525883
$ws.Range("M107").Value = 1562.54544
$ws.Range("N107").Value = -2529723
$ws.Range("H122").Value = 4501
$ws.Range("I122").Value = 3124.25
$ws.Range("J122").Value = 10008
$ws.Range("K122").Value = 9372.75
$ws.Range("L122").Value = 30024
$ws.Range("M122").Value = -6922.75
$ws.Range("N122").Value = -34924
$ws.Range("H132").Value = 3015.5862
$ws.Range("I132").Value = 2219
$ws.Range("J132").Value = 4144.0835
$ws.Range("K132").Value = 6657
$ws.Range("L132").Value = 12432.2505
$ws.Range("M132").Value = -4127
$ws.Range("H134").Value = 14460.2
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 14460.2
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 43380.60000000001
$ws.Range("N134").Value = -48450.60000000001

$ws = $wb.Worksheets("LTW")
$ws.Range("H16").Value = 15113389
$ws.Range("I16").Value = 42001000
$ws.Range("J16").Value = 1669583.5
$ws.Range("K16").Value = 42001000
$ws.Range("L16").Value = 1669583.5
$ws.Range("M16").Value = -42000830
$ws.Range("N16").Value = -1669923.5
$ws.Range("H55").Value = 253228.02
$ws.Range("I55").Value = 542045
$ws.Range("J55").Value = 513.1667
$ws.Range("K55").Value = 542045
$ws.Range("L55").Value = 513.1667
$ws.Range("M55").Value = -541872
$ws.Range("N55").Value = -859.1667
$ws.Range("H82").Value = 1619
$ws.Range("I82").Value = 1251.5454
$ws.Range("J82").Value = 2629.5
$ws.Range("K82").Value = 1251.5454
$ws.Range("L82").Value = 2629.5
$ws.Range("M82").Value = -890.5454
$ws.Range("N82").Value = -3351.5
$ws.Range("H85").Value = 1619
$ws.Range("I85").Value = 1251.5454
$ws.Range("J85").Value = 2629.5
$ws.Range("K85").Value = 1251.5454
$ws.Range("L85").Value = 2629.5
$ws.Range("M85").Value = -3.545399999999972
$ws.Range("N85").Value = -5125.5
$ws.Range("H122").Value = 3341.1853
$ws.Range("I122").Value = 3240.6
$ws.Range("J122").Value = 3628.5715
$ws.Range("K122").Value = 9721.8
$ws.Range("L122").Value = 10885.7145
$ws.Range("M122").Value = -7271.799999999999
$ws.Range("N122").Value = -15785.7145

$ws = $wb.Worksheets("WVR")
$ws.Range("H122").Value = 1959.8235
$ws.Range("I122").Value = 2064
$ws.Range("J122").Value = 1621.25
$ws.Range("K122").Value = 6192
$ws.Range("L122").Value = 4863.75
$ws.Range("M122").Value = -3742
$ws.Range("N122").Value = -9763.75
$ws.Range("H132").Value = 2972.7368
$ws.Range("I132").Value = 1494.7413
$ws.Range("J132").Value = 7735.1665
$ws.Range("K132").Value = 4484.2239
$ws.Range("L132").Value = 23205.4995
$ws.Range("M132").Value = -1954.2239

